$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "Info"
$ws.Range("A1").Value = "source: https://ae-scenario-explorer.cloud.set.kuleuven.be"
$ws.Move($null, $wb.Worksheets.Item("scenario"))
